# Weekly update: a new week's record is inserted at row 37. This pushes
# every previously-existing row from 37 downward down by one row (old
# row 148 ends up at row 149), growing the table by one row overall.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new row at position 37; Excel shifts rows 37..148 down to 38..149
# and grows the sheet dimension accordingly (A1:T148 -> A1:T149).
$ws.Rows("37:37").Insert()

# Populate the newly-inserted row 37 with the new week's data. Columns that
# are constant for every record in this sheet (A,B,C,E,F,G,H,I,J,R) keep the
# same values as the rest of the table.
$ws.Range("A37").Value = 11
$ws.Range("B37").Value = "Vega Monumental Concepción"
$ws.Range("C37").Value = "Bíobío"
$ws.Range("D37").Value = 44623
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100108
$ws.Range("H37").Value = "Tropicales y subtropicales"
$ws.Range("I37").Value = 100108005
$ws.Range("J37").Value = "Piña"
$ws.Range("K37").Value = "Caramelo"
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 250
$ws.Range("N37").Value = 15000
$ws.Range("O37").Value = 16000
$ws.Range("P37").Value = 15480
$ws.Range("Q37").Value = "$/caja 14 unidades"
$ws.Range("R37").Value = "Ecuador"
$ws.Range("S37").Value = 1106
$ws.Range("T37").Value = 14
